$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.803.35'
$ws.Range("E2").Value = '  +0.19%  '
$ws.Range("D3").Value = '3.164.34'
$ws.Range("E3").Value = '  +0.00%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '616.60'
$ws.Range("E5").Value = '  +2.53%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.20'
$ws.Range("E6").Value = '  -1.76%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").Value = '3.155.88'
$ws.Range("E8").Value = '  -0.29%  '
$ws.Range("E9").Value = '  -1.02%  '
$ws.Range("E10").Value = '  -1.12%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.47'
$ws.Range("E11").Value = '  -3.01%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.473'
$ws.Range("E12").Value = '  -1.36%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000259'
$ws.Range("E13").Value = '  -0.40%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.84'
$ws.Range("E14").Value = '  -3.55%  '
$ws.Range("D15").Value = '3.683.48'
$ws.Range("E15").Value = '  +0.14%  '
$ws.Range("E16").Value = '  +2.83%  '
$ws.Range("D17").Value = '64.786.12'
$ws.Range("E17").Value = '  +0.07%  '
$ws.Range("D18").Value = '3.162.61'
$ws.Range("E18").Value = '  -0.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.94'
$ws.Range("E19").Value = '  -1.78%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '481.79'
$ws.Range("E20").Value = '  -0.52%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.75'
$ws.Range("E21").Value = '  -1.04%  '
$ws.Range("E22").Value = '  +0.35%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.05'
$ws.Range("E23").Value = '  +2.91%  '
$ws.Range("E24").Value = '  -1.52%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.39'
$ws.Range("E25").Value = '  -1.01%  '
$ws.Range("E26").Value = '  +0.08%  '
$ws.Range("E27").Value = '  -2.94%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.58'
$ws.Range("E28").Value = '  -2.31%  '
$ws.Range("E29").Value = '  -2.33%  '
$ws.Range("E30").Value = '  -2.33%  '
$ws.Range("E31").Value = '  -7.88%  '
$ws.Range("B32").Value = 'FirstDigitalUSD'
$ws.Range("C32").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.00'
$ws.Range("E32").Value = '  +0.18%  '
$ws.Range("B33").Value = 'Stacks'
$ws.Range("C33").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.73'
$ws.Range("E33").Value = '  -0.54%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.56'
$ws.Range("E34").Value = '  -1.44%  '
$ws.Range("E35").Value = '  +1.91%  '
$ws.Range("D36").Value = '0.0₃0783'
$ws.Range("E36").Value = '  +4.64%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.02'
$ws.Range("E37").Value = '  -2.32%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.24'
$ws.Range("E38").Value = '  -0.30%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '53.08'
$ws.Range("E39").Value = '  -3.24%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '461.77'
$ws.Range("E40").Value = '  +0.28%  '
$ws.Range("E41").Value = '  -0.65%  '
$ws.Range("E42").Value = '  -4.54%  '
$ws.Range("E43").Value = '  -1.68%  '
$ws.Range("D44").Value = '2.847.73'
$ws.Range("E44").Value = '  -1.77%  '
$ws.Range("E45").Value = '  -4.38%  '
$ws.Range("E46").Value = '  -2.60%  '
$ws.Range("E47").Value = '  +4.89%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '26.68'
$ws.Range("E48").Value = '  -1.71%  '
$ws.Range("E49").Value = '  +0.10%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.114'
$ws.Range("E50").Value = '  -1.49%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '120.62'
$ws.Range("E51").Value = '  +0.73%  '
